$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "MG Contécnica - Unidade Barra da Tijuca"
$ws.Range("B2").Value = ""

# Row 3
$ws.Range("A3").Value = "AXM"
$ws.Range("B3").Value = "Tv. do Ouvidor, 5 - 4 andar - Centro, Rio de Janeiro - RJ, 20040-040"
$ws.Range("C3").Value = "(21) 2206-1000"
$ws.Range("D3").Value = "axms.com.br"

# Row 4
$ws.Range("A4").Value = "Seletus Contabilidade - RJ"
$ws.Range("B4").Value = "R. Campo Grande, 1014 - sala 526 - Campo Grande, Rio de Janeiro - RJ, 23080-000"
$ws.Range("C4").Value = "(21) 4107-1417"
$ws.Range("D4").Value = "seletuscontabilidade.com.br"

# Row 5
$ws.Range("A5").Value = "Seletus Contabilidade - RJ"
$ws.Range("B5").Value = "R. Campo Grande, 1014 - sala 526 - Campo Grande, Rio de Janeiro - RJ, 23080-000"
$ws.Range("C5").Value = "(21) 4107-1417"
$ws.Range("D5").Value = "seletuscontabilidade.com.br"

# Row 6
$ws.Range("A6").Value = "Cemage Contabilidade"
$ws.Range("B6").Value = "R. Viúva Dantas, 60 - Campo Grande, Rio de Janeiro - RJ, 23050-090"
$ws.Range("C6").Value = "(21) 2413-5334"
